$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 50
$wsSummary.Range("E4").Value = 50
# Update the sheet's remembered selection (without permanently changing the active tab)
$wsSummary.Range("A7:XFD14").Select()

# --- Repayment schedule sheet updates ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("I5").Value = 0
$wsRepay.Range("K5").Value = 887.72
$wsRepay.Range("P5").Value = 887.72
# Update the sheet's remembered selection (without permanently changing the active tab)
$wsRepay.Range("A9:XFD9").Select()

# Restore the originally active sheet/tab (Transactions) so the workbook-level
# active tab and tabSelected flags are left untouched by the selection updates above.
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()
